# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" sheet before "ODI Batting" with player bio data.
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE on both ODI Batting and ODI Bowling
#    sheets, replacing the full scorecard URL with just the numeric match code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the "Player Info" sheet in front of "ODI Batting"
#    (worksheet references resolve positionally, so re-fetch sheets by name
#    after every operation that inserts/removes/reorders sheets)
# ---------------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$infoSheet.Name = "Player Info"

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $infoHeaders.Length; $col++) {
    $cell = $infoSheet.Cells.Item(1, $col)
    $cell.Value = $infoHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$infoValues = @("4647", "Glenn Dominic Phillips", "Right Handed", "Does Not Bowl | Unknown")
for ($col = 1; $col -le $infoValues.Length; $col++) {
    $cell = $infoSheet.Cells.Item(2, $col)
    # Force text storage so numeric-looking values (e.g. the ID) stay strings.
    $cell.NumberFormat = "@"
    $cell.Value = $infoValues[$col - 1]
}

$infoSheet.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Replace MATCH_CARD_LINK with MATCH_CODE on "ODI Batting" (column D)
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

for ($row = 2; $row -le 17; $row++) {
    $cell = $battingSheet.Cells.Item($row, 4)
    $link = $cell.Value2
    if ($link -ne $null -and $link -ne "") {
        $code = $link -replace '.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# ---------------------------------------------------------------------------
# 3. Replace MATCH_CARD_LINK with MATCH_CODE on "ODI Bowling" (column B)
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

for ($row = 2; $row -le 7; $row++) {
    $cell = $bowlingSheet.Cells.Item($row, 2)
    $link = $cell.Value2
    if ($link -ne $null -and $link -ne "") {
        $code = $link -replace '.*MatchCode=', ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}
